# VerifyUINoCFCCTestData.xlsx - "Added BWP SelIDE Test Cases, Fixed VRelay 2.5 Issues"
#
# The CC* (credit-card) prefixed headers used by the "CC" verification flow
# are renamed to their plain equivalents (CCAmount -> Amount, CCUDF1 -> UDF1,
# ... CCEmail -> Email) on both sheets. CCDate is intentionally left as-is.
# Sheet2 becomes the active/selected sheet, selections move to the header
# row, and a "Public" classification footer is stamped onto both sheets.

$wb = $excel.ActiveWorkbook

$headerRenames = @{
    "B1" = "Amount"
    "C1" = "UDF1"
    "D1" = "UDF2"
    "E1" = "UDF3"
    "F1" = "UDF4"
    "G1" = "UDF5"
    "H1" = "UDF6"
    "I1" = "UDF7"
    "J1" = "UDF8"
    "K1" = "UDF9"
    "L1" = "UDF10"
    "M1" = "Name"
    "N1" = "CardNum"
    "O1" = "SPC"
    "P1" = "ExpM"
    "Q1" = "ExpY"
    "S1" = "AL1"
    "T1" = "AL2"
    "U1" = "ZIP"
    "V1" = "Email"
}

# Classification footer stamped onto both sheets (carriage-return + font/size
# / colour codes ahead of the literal "Public" marker).
$cr = [char]13
$footerText = "$cr&1#&`"Calibri`"&10&K000000 Public "

for ($sheetIdx = 1; $sheetIdx -le 2; $sheetIdx++) {
    $ws = $wb.Worksheets.Item($sheetIdx)

    foreach ($addr in $headerRenames.Keys) {
        $ws.Range($addr).Value = $headerRenames[$addr]
    }

    $ws.PageSetup.CenterFooter = $footerText
}

# --- Sheet1: drop the frozen-pane tab selection, move the plain selection
#     onto the header row ------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Activate()
$window = $excel.ActiveWindow
$window.FreezePanes = $false
$ws1.Range("B1").Select()
$window.FreezePanes = $true
$ws1.Range("B1:V1").Select()

# --- Sheet2: becomes the active/selected tab, selection moves onto the
#     header row too ------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Activate()
$ws2.Range("B1:V1").Select()
